$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hyperparameter Optimization")

# Populate cells that introduce brand-new shared strings first, in the same
# order the source workbook authored them, so the shared-string table layout matches.
$ws.Range("P53").Value = "[0.06173151358962059, 0.046613264828920364, 0.07495468109846115, 0.08225483447313309, 0.06503096967935562]"
$ws.Range("K53").Value = "40"
$ws.Range("K13").Value = "48"
$ws.Range("P13").Value = "[0.06516171991825104, 0.05991170182824135, 0.08209428936243057, 0.0976945087313652, 0.08277694880962372]"
$ws.Range("P54").Value = "[0.047666534781455994, 0.05321182310581207, 0.06995537877082825, 0.08843991160392761, 0.0692932978272438]"
$ws.Range("K56").Value = "43"
$ws.Range("P56").Value = "[0.042313531041145325, 0.05110793933272362, 0.06351163238286972, 0.06419646739959717, 0.05570153892040253]"
$ws.Range("K16").Value = "55"
$ws.Range("P16").Value = "[0.05161529406905174, 0.04834244027733803, 0.07278178632259369, 0.08038552850484848, 0.06707189232110977]"
$ws.Range("K14").Value = "4"
$ws.Range("P14").Value = "[0.06150340661406517, 0.060546960681676865, 0.08885558694601059, 0.10247880965471268, 0.07931020110845566]"

# Remaining cells in each row (numbers + already-existing text values)
# Row 13
$ws.Range("F13").Value = 0.30861173110805101
$ws.Range("G13").Value = "Adam"
$ws.Range("H13").Value = "relu"
$ws.Range("I13").Value = "1024"
$ws.Range("J13").Value = 0.44613235732750101
$ws.Range("L13").Value = 0.077527833729982301
$ws.Range("M13").Value = 0.0135524859768093

# Row 14
$ws.Range("F14").Value = 0.0086635758205824993
$ws.Range("G14").Value = "SGD"
$ws.Range("H14").Value = "sigmoid"
$ws.Range("I14").Value = "1024"
$ws.Range("J14").Value = 0.29873396494652799
$ws.Range("L14").Value = 0.078538993000984103
$ws.Range("M14").Value = 0.016087673072951801

# Row 16
$ws.Range("F16").Value = 0.0220349393162898
$ws.Range("G16").Value = "Adam"
$ws.Range("H16").Value = "relu"
$ws.Range("I16").Value = "256"
$ws.Range("J16").Value = 0.48093551450042699
$ws.Range("L16").Value = 0.064039388298988295
$ws.Range("M16").Value = 0.0122765972329933

# Row 53
$ws.Range("E53").Value = "512"
$ws.Range("F53").Value = 0.041375766681849703
$ws.Range("G53").Value = "SGD"
$ws.Range("H53").Value = "tanh"
$ws.Range("I53").Value = "1024"
$ws.Range("J53").Value = 0.156686443776829
$ws.Range("L53").Value = 0.066117052733898105
$ws.Range("M53").Value = 0.012160114696852299

# Row 54
$ws.Range("E54").Value = "1024"
$ws.Range("F54").Value = 0.071322500106832706
$ws.Range("G54").Value = "SGD"
$ws.Range("H54").Value = "relu"
$ws.Range("I54").Value = "512"
$ws.Range("J54").Value = 0.16079535852270399
$ws.Range("K54").Value = "92"
$ws.Range("L54").Value = 0.065713389217853496
$ws.Range("M54").Value = 0.0143477047020499

# Row 56
$ws.Range("E56").Value = "512"
$ws.Range("F56").Value = 0.055894697755093099
$ws.Range("G56").Value = "SGD"
$ws.Range("H56").Value = "tanh"
$ws.Range("I56").Value = "512"
$ws.Range("J56").Value = 0.273986074649267
$ws.Range("L56").Value = 0.055366221815347602
$ws.Range("M56").Value = 0.0081601359559920292


# Update the saved view position/selection to match the latest edit session
$ws.Activate()
$ws.Range("L15").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
